# LOG-T28-Verify validation for invalid email format (missing @)
#
# Sheet1: the "Password"/"ConfirmPassword" test rows are updated so the
# Password column (D) no longer mirrors the ConfirmPassword column (E) -
# new password-like values are entered. Two of the three contain an "@",
# so Excel auto-recognizes them as email-like text and turns them into
# mailto hyperlinks (Hyperlink cell style). The DuplicateEmail sheet gets
# the same treatment for its two rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet1
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("D2").Value = "Ravi@2025"
$ws1.Hyperlinks.Add($ws1.Range("D2"), "mailto:Ravi@2025")
$ws1.Range("D2").Style = "Hyperlink"

$ws1.Range("D3").Value = "Maria#3215"

$ws1.Range("D4").Value = "Ravi@20245"
$ws1.Hyperlinks.Add($ws1.Range("D4"), "mailto:Ravi@20245")
$ws1.Range("D4").Style = "Hyperlink"

# ---------------------------------------------------------------
# DuplicateEmail sheet
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("DuplicateEmail")

$ws3.Range("D2").Value = "Ravi@2025"
$ws3.Hyperlinks.Add($ws3.Range("D2"), "mailto:Ravi@2025")
$ws3.Range("D2").Style = "Hyperlink"

$ws3.Range("D3").Value = "Ravi@2025"
$ws3.Hyperlinks.Add($ws3.Range("D3"), "mailto:Ravi@2025")
$ws3.Range("D3").Style = "Hyperlink"

# Selection on DuplicateEmail moves to I9 (set before re-activating Sheet1
# below so Sheet1 ends up as the tab that is actually selected).
$ws3.Range("I9").Select()

# ---------------------------------------------------------------
# Final view state: Sheet1 active, selection on D6
# ---------------------------------------------------------------
$ws1.Activate()
$ws1.Range("D6").Select()
